$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ============================================================================
# "aggiunta 4 e 5 drupa 2021" - add the 2021 harvest ("drupa") sampling rows,
# rows 4 and 5 (BN codes) included, mirroring the 2020 block; also renames the
# old BN_1..BN_5 codice_azienda codes to the new zero-padded BN_01..BN_05 form.
# ============================================================================

# --- Step 1: rename the 2020 "BN_1".."BN_5" codes to "BN_01".."BN_05" ---
$ws.Range("A8").Value = "BN_01"
$ws.Range("A9").Value = "BN_02"
$ws.Range("A10").Value = "BN_03"
$ws.Range("A11").Value = "BN_04"
$ws.Range("A12").Value = "BN_05"

# --- Step 2: append the 2021 rows 28-53 (values only; styles applied below) ---

# Row 28
$ws.Range("A28").Value = "AV_01"
$ws.Range("B28").Value = 2021
$ws.Range("C28").Value = "R2"
$ws.Range("D28").Value = 44551
$ws.Range("E28").Value = "SI"
$ws.Range("F28").Value = "NO"

# Row 29
$ws.Range("A29").Value = "AV_02"
$ws.Range("B29").Value = 2021
$ws.Range("C29").Value = "R2"
$ws.Range("D29").Value = 44551
$ws.Range("E29").Value = "SI"
$ws.Range("F29").Value = "NO"

# Row 30
$ws.Range("A30").Value = "AV_03"
$ws.Range("B30").Value = 2021
$ws.Range("C30").Value = "R2"
$ws.Range("D30").Value = 44551
$ws.Range("E30").Value = "SI"
$ws.Range("F30").Value = "NO"

# Row 31
$ws.Range("A31").Value = "AV_04"
$ws.Range("B31").Value = 2021
$ws.Range("C31").Value = "R2"
$ws.Range("D31").Value = 44551
$ws.Range("E31").Value = "SI"
$ws.Range("F31").Value = "NO"

# Row 32 (highlighted red)
$ws.Range("A32").Value = "AV_05"
$ws.Range("B32").Value = 2021
$ws.Range("C32").Value = "R2"

# Row 33
$ws.Range("A33").Value = "AV_06"
$ws.Range("B33").Value = 2021
$ws.Range("C33").Value = "R2"
$ws.Range("D33").Value = 44551
$ws.Range("E33").Value = "SI"
$ws.Range("F33").Value = "NO"

# Row 34
$ws.Range("A34").Value = "BN_01"
$ws.Range("B34").Value = 2021
$ws.Range("C34").Value = "R2"
$ws.Range("D34").Value = 44523
$ws.Range("E34").Value = "SI"
$ws.Range("F34").Value = "SI"

# Row 35
$ws.Range("A35").Value = "BN_02"
$ws.Range("B35").Value = 2021
$ws.Range("C35").Value = "R2"
$ws.Range("D35").Value = 44538
$ws.Range("E35").Value = "SI"
$ws.Range("F35").Value = "SI"

# Row 36 (highlighted red)
$ws.Range("A36").Value = "BN_03"
$ws.Range("B36").Value = 2021
$ws.Range("C36").Value = "R2"

# Row 37
$ws.Range("A37").Value = "BN_04"
$ws.Range("B37").Value = 2021
$ws.Range("C37").Value = "R2"
$ws.Range("D37").Value = 44516
$ws.Range("E37").Value = "SI"
$ws.Range("F37").Value = "SI"

# Row 38
$ws.Range("A38").Value = "BN_05"
$ws.Range("B38").Value = 2021
$ws.Range("C38").Value = "R2"
$ws.Range("D38").Value = 44531
$ws.Range("E38").Value = "SI"
$ws.Range("F38").Value = "SI"

# Row 39
$ws.Range("A39").Value = "CE_01"
$ws.Range("B39").Value = 2021
$ws.Range("C39").Value = "R2"
$ws.Range("D39").Value = 44549
$ws.Range("E39").Value = "SI"
$ws.Range("F39").Value = "NO"

# Row 40 (highlighted red)
$ws.Range("A40").Value = "CE_02"
$ws.Range("B40").Value = 2021
$ws.Range("C40").Value = "R2"

# Row 41
$ws.Range("A41").Value = "CE_03"
$ws.Range("B41").Value = 2021
$ws.Range("C41").Value = "R2"

# Row 42
$ws.Range("A42").Value = "NA_01"
$ws.Range("B42").Value = 2021
$ws.Range("C42").Value = "R2"

# Row 43
$ws.Range("A43").Value = "NA_02"
$ws.Range("B43").Value = 2021
$ws.Range("C43").Value = "R2"

# Row 44
$ws.Range("A44").Value = "SA_01"
$ws.Range("B44").Value = 2021
$ws.Range("C44").Value = "R2"
$ws.Range("D44").Value = 44546
$ws.Range("E44").Value = "SI"
$ws.Range("F44").Value = "NO"

# Row 45
$ws.Range("A45").Value = "SA_02"
$ws.Range("B45").Value = 2021
$ws.Range("C45").Value = "R2"
$ws.Range("D45").Value = 44546
$ws.Range("E45").Value = "SI"
$ws.Range("F45").Value = "NO"

# Row 46
$ws.Range("A46").Value = "SA_03"
$ws.Range("B46").Value = 2021
$ws.Range("C46").Value = "R2"
$ws.Range("D46").Value = 44557
$ws.Range("E46").Value = "SI"
$ws.Range("F46").Value = "NO"

# Row 47
$ws.Range("A47").Value = "SA_04"
$ws.Range("B47").Value = 2021
$ws.Range("C47").Value = "R2"
$ws.Range("D47").Value = 44552
$ws.Range("E47").Value = "SI"
$ws.Range("F47").Value = "NO"

# Row 48
$ws.Range("A48").Value = "SA_05"
$ws.Range("B48").Value = 2021
$ws.Range("C48").Value = "R2"

# Row 49
$ws.Range("A49").Value = "SA_06"
$ws.Range("B49").Value = 2021
$ws.Range("C49").Value = "R2"

# Row 50
$ws.Range("A50").Value = "SA_07"
$ws.Range("B50").Value = 2021
$ws.Range("C50").Value = "R2"
$ws.Range("D50").Value = 44550
$ws.Range("E50").Value = "SI"
$ws.Range("F50").Value = "NO"

# Row 51
$ws.Range("A51").Value = "SA_08"
$ws.Range("B51").Value = 2021
$ws.Range("C51").Value = "R2"
$ws.Range("D51").Value = 44557
$ws.Range("E51").Value = "SI"
$ws.Range("F51").Value = "NO"

# Row 52 (highlighted red)
$ws.Range("A52").Value = "SA_09"
$ws.Range("B52").Value = 2021
$ws.Range("C52").Value = "R2"

# Row 53
$ws.Range("A53").Value = "SA_10"
$ws.Range("B53").Value = 2021
$ws.Range("C53").Value = "R2"

# --- Step 3: G28/G29/G30 olive-oil descriptions must be written in this exact
#     order (blend, then Ravece, then Marinese) so the new shared strings are
#     appended in the same order as in the target workbook ---
$ws.Range("G29").Value = "Olio extra vergine di oliva filtrato blend"
$ws.Range("G28").Value = "Olio extra vergine di oliva filtrato monovarietale Ravece"
$ws.Range("G30").Value = "Olio extra vergine di oliva filtrato monovarietale Marinese"
$ws.Range("G31").Value = "Olio extra vergine di oliva filtrato monovarietale Ravece"
$ws.Range("G33").Value = "Olio extra vergine di oliva filtrato monovarietale Ravece"

# --- Step 4: give every new date cell in column D the same date style already
#     used by D2:D27 (number format 14), without creating a duplicate numFmt ---
$dateTargets = $excel.Union($ws.Range("D28"), $ws.Range("D29"), $ws.Range("D30"), $ws.Range("D31"), $ws.Range("D33"), $ws.Range("D34"), $ws.Range("D35"), $ws.Range("D37"), $ws.Range("D38"), $ws.Range("D39"), $ws.Range("D44"), $ws.Range("D45"), $ws.Range("D46"), $ws.Range("D47"), $ws.Range("D50"), $ws.Range("D51"))
$ws.Range("D2").Copy()
$dateTargets.PasteSpecial(-4122)

# --- Step 5: highlight the 4 special rows (no date recorded yet) with the same
#     red fill used elsewhere in the workbook ---
$highlightTargets = $excel.Union($ws.Range("A32:C32"), $ws.Range("A36:C36"), $ws.Range("A40:C40"), $ws.Range("A52:C52"))
$highlightTargets.Interior.Color = 255

# --- Step 6: clear the clipboard marquee and restore the expected selection ---
$excel.CutCopyMode = 0
$ws.Range("F37").Select()
